$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 243
$ws.Range("C2").Value = 84.08
$ws.Range("B3").Value = 46
$ws.Range("C3").Value = 15.92
